# Update "想去人数" (F column) values across sheets, as produced by the
# site's generator run (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 1387
$wsExhibit.Range("F7").Value = 141
$wsExhibit.Range("F8").Value = 28
$wsExhibit.Range("F10").Value = 9415
$wsExhibit.Range("F12").Value = 100
$wsExhibit.Range("F13").Value = 220
$wsExhibit.Range("F15").Value = 358
$wsExhibit.Range("F16").Value = 6397
$wsExhibit.Range("F18").Value = 95
$wsExhibit.Range("F20").Value = 141

# --- 演出 (Performance) sheet ---
$wsPerform = $wb.Worksheets.Item("演出")
$wsPerform.Range("F2").Value = 37

# --- 全部类型 (All Types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 1387
$wsAll.Range("F7").Value = 141
$wsAll.Range("F8").Value = 28
$wsAll.Range("F10").Value = 37
$wsAll.Range("F12").Value = 9415
$wsAll.Range("F14").Value = 100
$wsAll.Range("F15").Value = 220
$wsAll.Range("F17").Value = 358
$wsAll.Range("F18").Value = 6397
$wsAll.Range("F20").Value = 95
$wsAll.Range("F22").Value = 141
